$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above row 2 (old Dutatama row) for the new students
$ws.Rows("2:7").Insert()
$ws.Rows("2:7").ClearFormats()

# Fill in the new student rows (2-7), all class 3SD2 — entered column by column
$names   = @("Andika", "Sabila", "Fauzan", "Ajeng", "Khesya", "Brigitta")
$nims    = @(999888777, 888777666, 777666555, 666555444, 555444333, 444333222)
$emails  = @("999888777@stis.ac.id", "888777666@stis.ac.id", "777666555@stis.ac.id", "666555444@stis.ac.id", "555444333@stis.ac.id", "444333222@stis.ac.id")
$genders = @("Laki-laki", "Perempuan", "Laki-laki", "Perempuan", "Perempuan", "Perempuan")
$classes = @("3SD2", "3SD2", "3SD2", "3SD2", "3SD2", "3SD2")

for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 2, 1).Value = $names[$i] }
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 2, 2).Value = $nims[$i] }
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 2, 3).Value = $emails[$i] }
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 2, 4).Value = $genders[$i] }
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 2, 5).Value = $classes[$i] }

# Re-type (fix) the email column for the old rows (now shifted to rows 8-13)
$oldEmails = @("111222333@stis.ac.id", "222333444@stis.ac.id", "333444555@stis.ac.id", "444555666@stis.ac.id", "666777888@stis.ac.id", "777888999@stis.ac.id")
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 8, 3).ClearContents() }
for ($i = 0; $i -lt 6; $i++) { $ws.Cells.Item($i + 8, 3).Value = $oldEmails[$i] }

$ws.Range("C5").Select()
